$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2023-09-04 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-05 Tuesday", 2) | Out-Null

# Update each table cell value (row-major order matches the diff)
$t = $d.Tables.Item(1)
$newValues = @(
  @("34+19=53", "29+55=84", "24+40=64", "8+14=22", "32-18=14"),
  @("38+28=66", "33+51=84", "51-23=28", "82-42=40", "47+0=47"),
  @("18-13=5", "44+2=46", "36-31=5", "61-6=55", "80+8=88"),
  @("5+20=25", "98-44=54", "71-62=9", "9+90=99", "39+22=61"),
  @("8+51=59", "31+28=59", "54+12=66", "38+8=46", "71-16=55"),
  @("39+49=88", "33+3=36", "15+73=88", "60+9=69", "51-24=27"),
  @("69+24=93", "58-30=28", "62-40=22", "32+35=67", "89-81=8"),
  @("72-43=29", "48+30=78", "44-36=8", "21+34=55", "59+25=84"),
  @("81-31=50", "66-15=51", "55+5=60", "54-39=15", "72+18=90"),
  @("73-69=4", "14+61=75", "25+27=52", "56-27=29", "68-25=43"),
  @("44+0=44", "29+13=42", "29+33=62", "80-16=64", "64-4=60"),
  @("11+56=67", "53-10=43", "0+5=5", "10+43=53", "27+50=77"),
  @("53-15=38", "24+52=76", "5+53=58", "80-36=44", "17+69=86"),
  @("4+46=50", "17+37=54", "36+0=36", "89-58=31", "74-33=41"),
  @("44+37=81", "58+7=65", "24+71=95", "5+44=49", "60-18=42"),
  @("32+2=34", "2+23=25", "26+40=66", "97-35=62", "2+40=42"),
  @("19-2=17", "78+16=94", "54-47=7", "56+35=91", "30+64=94"),
  @("30+42=72", "50-37=13", "73-65=8", "87-85=2", "46+32=78"),
  @("49-32=17", "27+50=77", "47+13=60", "16+22=38", "4+87=91"),
  @("26+60=86", "54-36=18", "14+11=25", "74-55=19", "0+90=90")
)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le 5; $c++) {
    $t.Rows.Item($r).Cells.Item($c).Range.Text = $newValues[$r-1][$c-1]
  }
}
